{"js": "// The duration in the first paragraph changes from \"3 months\" to \"4 months\",\n// and the Word-managed \"_GoBack\" bookmark (last-edit marker) moves from its\n// old spot (right after \"session (\") to the new edit location (right after\n// the \"4\" that replaced the \"3\").\n\n// 1) Drop the stale \"_GoBack\" bookmark that currently sits near \"session (\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Replace the \"3\" in \"3 months\" (first paragraph) with \"4\".\nconst firstParagraph = context.document.body.paragraphs.getFirst();\nconst threeResults = firstParagraph.search(\"3\", { matchCase: true, matchWholeWord: true });\nthreeResults.load(\"items\");\nawait context.sync();\n\nconst threeRange = threeResults.items[0];\nthreeRange.insertText(\"4\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Re-find the \"4\" we just inserted and drop a fresh \"_GoBack\" bookmark\n//    right after it, so the paragraph ends up split into \"4\" + \" months\"\n//    runs with the bookmark in between.\nconst fourResults = context.document.body.paragraphs.getFirst()\n  .search(\"4\", { matchCase: true, matchWholeWord: true });\nfourResults.load(\"items\");\nawait context.sync();\n\nconst fourRange = fourResults.items[0];\nconst afterFour = fourRange.getRange(Word.RangeLocation.end);\nafterFour.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The duration in the first paragraph changes from \"3 months\" to \"4 months\".\n# Word's auto-managed \"_GoBack\" bookmark (marks the last edit location) needs\n# to move from its old spot (right after \"session (\") to the new edit spot\n# (right after the \"4\" that replaced the \"3\").\n\n$d = $word.ActiveDocument\n\n# 1) Remove the stale \"_GoBack\" bookmark currently sitting near \"session (\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Replace the \"3\" in \"3 months\" (the document's first paragraph) with \"4\".\n$firstPara = $d.Paragraphs.Item(1)\n$rng = $firstPara.Range\n$find = $rng.Find\n$find.Text = \"3\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n[void]$find.Execute()\n$rng.Text = \"4\"\n\n# 3) Collapse the range to right after the \"4\" and drop a fresh \"_GoBack\"\n#    bookmark there, so the paragraph ends up split into \"4\" + \" months\"\n#    runs with the bookmark in between.\n$rng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n"}
